$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: was ikdossantoscs@gmail.com / ikram2 / ikram123
#        -> retnowardani024@gmail.com / sekwan / gaul
$ws.Range("A2").Value = "retnowardani024@gmail.com"
$ws.Range("B2").Value = "sekwan"
$ws.Range("C2").Value = "gaul"

# Row 3: was retnowardani024@gmail.com / retno gaul / admin34
#        -> ikdossantoscs@gmail.com / ikram / ikram123
$ws.Range("A3").Value = "ikdossantoscs@gmail.com"
$ws.Range("B3").Value = "ikram"
$ws.Range("C3").Value = "ikram123"
